$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full replacement data for rows 2-17 (columns A-T), reflecting the addition
# of a new "M2" sending/target cluster group to the Ntn4-Unc5a LR-pair table.
$data = @(
    @("ECs", "Ntn4", "Unc5a", "ECs", 3, 1, 2.75977, 8.279310000000001, 0.06698501710740076, 0.06698501710740076, 3, 1, 0.837459, 2.512377, 0.1861007731134217, 0.1861007731134217, 2.31119422443, 20.80074801987, 0.01246596347070306, 0.01246596347070306),
    @("ECs", "Ntn4", "Unc5a", "FAPs", 3, 1, 2.75977, 8.279310000000001, 0.06698501710740076, 0.06698501710740076, 3, 1, 1.864510666666667, 5.593532, 0.4143329721752205, 0.4143329721752205, 5.145620602546667, 46.31058542292, 0.02775410122931735, 0.02775410122931735),
    @("ECs", "Ntn4", "Unc5a", "M2", 3, 1, 2.75977, 8.279310000000001, 0.06698501710740076, 0.06698501710740076, 3, 1, 0.4771786666666666, 1.431536, 0.1060390046317472, 0.1060390046317472, 1.316903368906667, 11.85213032016, 0.007103024539309338, 0.007103024539309338),
    @("ECs", "Ntn4", "Unc5a", "sCs", 3, 1, 2.75977, 8.279310000000001, 0.06698501710740076, 0.06698501710740076, 3, 1, 1.320881333333333, 3.962644, 0.2935272500796106, 0.2935272500796106, 3.645328677293334, 32.80795809564, 0.01966192786807102, 0.01966192786807102),
    @("FAPs", "Ntn4", "Unc5a", "ECs", 3, 1, 18.996877, 56.99063100000001, 0.4610913702345443, 0.4610913702345442, 3, 1, 0.837459, 2.512377, 0.1861007731134217, 0.1861007731134217, 15.909105615543, 143.181950539887, 0.08580946047657564, 0.08580946047657563),
    @("FAPs", "Ntn4", "Unc5a", "FAPs", 3, 1, 18.996877, 56.99063100000001, 0.4610913702345443, 0.4610913702345442, 3, 1, 1.864510666666667, 5.593532, 0.4143329721752205, 0.4143329721752205, 35.41987979985467, 318.7789181986921, 0.1910453578736237, 0.1910453578736237),
    @("FAPs", "Ntn4", "Unc5a", "M2", 3, 1, 18.996877, 56.99063100000001, 0.4610913702345443, 0.4610913702345442, 3, 1, 0.4771786666666666, 1.431536, 0.1060390046317472, 0.1060390046317472, 9.064904437690666, 81.584139939216, 0.04889366994395952, 0.04889366994395952),
    @("FAPs", "Ntn4", "Unc5a", "sCs", 3, 1, 18.996877, 56.99063100000001, 0.4610913702345443, 0.4610913702345442, 3, 1, 1.320881333333333, 3.962644, 0.2935272500796106, 0.2935272500796106, 25.09262022092934, 225.833581988364, 0.1353428819403854, 0.1353428819403854),
    @("M2", "Ntn4", "Unc5a", "ECs", 1, 0.3333333333333333, 0.07467966666666666, 0.224039, 0.001812621613120533, 0.001812621613120533, 3, 1, 0.837459, 2.512377, 0.1861007731134217, 0.1861007731134217, 0.06254115896699999, 0.5628704307029999, 0.0003373302835638287, 0.0003373302835638287),
    @("M2", "Ntn4", "Unc5a", "FAPs", 1, 0.3333333333333333, 0.07467966666666666, 0.224039, 0.001812621613120533, 0.001812621613120533, 3, 1, 1.864510666666667, 5.593532, 0.4143329721752205, 0.4143329721752205, 0.1392410350831111, 1.253169315748, 0.000751028900393273, 0.0007510289003932729),
    @("M2", "Ntn4", "Unc5a", "M2", 1, 0.3333333333333333, 0.07467966666666666, 0.224039, 0.001812621613120533, 0.001812621613120533, 3, 1, 0.4771786666666666, 1.431536, 0.1060390046317472, 0.1060390046317472, 0.0356355437671111, 0.320719893904, 0.0001922085916292933, 0.0001922085916292933),
    @("M2", "Ntn4", "Unc5a", "sCs", 1, 0.3333333333333333, 0.07467966666666666, 0.224039, 0.001812621613120533, 0.001812621613120533, 3, 1, 1.320881333333333, 3.962644, 0.2935272500796106, 0.2935272500796106, 0.09864297767955554, 0.887786799116, 0.0005320538375341379, 0.0005320538375341377),
    @("sCs", "Ntn4", "Unc5a", "ECs", 3, 1, 19.36848366666667, 58.105451, 0.4701109910449345, 0.4701109910449345, 3, 1, 0.837459, 2.512377, 0.1861007731134217, 0.1861007731134217, 16.220310963003, 145.982798667027, 0.08748801888257918, 0.08748801888257916),
    @("sCs", "Ntn4", "Unc5a", "FAPs", 3, 1, 19.36848366666667, 58.105451, 0.4701109910449345, 0.4701109910449345, 3, 1, 1.864510666666667, 5.593532, 0.4143329721752205, 0.4143329721752205, 36.11274439365911, 325.014699542932, 0.1947824841718862, 0.1947824841718861),
    @("sCs", "Ntn4", "Unc5a", "M2", 3, 1, 19.36848366666667, 58.105451, 0.4701109910449345, 0.4701109910449345, 3, 1, 0.4771786666666666, 1.431536, 0.1060390046317472, 0.1060390046317472, 9.242227211415111, 83.180044902736, 0.0498501015568491, 0.0498501015568491),
    @("sCs", "Ntn4", "Unc5a", "sCs", 3, 1, 19.36848366666667, 58.105451, 0.4701109910449345, 0.4701109910449345, 3, 1, 1.320881333333333, 3.962644, 0.2935272500796106, 0.2935272500796106, 25.58346853027156, 230.251216772444, 0.1379903864336201, 0.1379903864336201)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
    }
}
